$wb = $excel.ActiveWorkbook

# The "Netherlands" sheet is the template the original author copied to
# create each new market sheet (same layout/styles as the other single-row
# market sheets, incl. the "bestFit" column D width and wrapped row 2).
$template = $wb.Worksheets.Item("Netherlands")

# --- Russia ---------------------------------------------------------------
$afterSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$template.Copy($null, $afterSheet)
$russia = $wb.Worksheets.Item($wb.Worksheets.Count)
$russia.Name = "Russia"
$russia.Range("B4").Value = "NGC-2929/T3320"
$russia.Range("B2").Value = "Russia Market"
$russia.Rows.Item(2).AutoFit()
$russia.Range("A1:D10").Select()

# --- Finland ----------------------------------------------------------------
$afterSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$template.Copy($null, $afterSheet)
$finland = $wb.Worksheets.Item($wb.Worksheets.Count)
$finland.Name = "Finland"
$finland.Range("B4").Value = "NGC-3130/T2957"
$finland.Range("B2").Value = "Finland Market"
$finland.Rows.Item(2).AutoFit()
$finland.Range("A1:D10").Select()

# --- Hungary ------------------------------------------------------------
$afterSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$template.Copy($null, $afterSheet)
$hungary = $wb.Worksheets.Item($wb.Worksheets.Count)
$hungary.Name = "Hungary"
$hungary.Range("B4").Value = "NGC-3104/T2992"
$hungary.Range("B2").Value = "Hungary Market"
$hungary.Rows.Item(2).AutoFit()

# Hungary ends up the active/selected sheet & tab, mirroring the committed
# workbook (firstSheet/activeTab point at the new last sheet).
$hungary.Activate()
$hungary.Range("J11").Select()
